$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 585.8214
$ws.Range("J17").Value = 585.8214
$ws.Range("L17").Value = 1757.4642
$ws.Range("N17").Value = -2093.4642
$ws.Range("H18").Value = 1232.8334
$ws.Range("I18").Value = 779.4
$ws.Range("K18").Value = 779.4
$ws.Range("M18").Value = -495.4
$ws.Range("H108").Value = 99906.664
$ws.Range("J108").Value = 99906.664
$ws.Range("L108").Value = 99906.664
$ws.Range("N108").Value = -107586.664
$ws.Range("H109").Value = 46890.777
$ws.Range("J109").Value = 46890.777
$ws.Range("L109").Value = 46890.777
$ws.Range("N109").Value = -49664.777
$ws.Range("H110").Value = 50196.168
$ws.Range("J110").Value = 50196.168
$ws.Range("L110").Value = 50196.168
$ws.Range("N110").Value = -58376.168
$ws.Range("H114").Value = 99741.664
$ws.Range("J114").Value = 99741.664
$ws.Range("L114").Value = 99741.664
$ws.Range("N114").Value = -108419.664
$ws.Range("H123").Value = 58818.625
$ws.Range("J123").Value = 58818.625
$ws.Range("L123").Value = 58818.625
$ws.Range("N123").Value = -68618.625
$ws.Range("H133").Value = 80871.375
$ws.Range("J133").Value = 80871.375
$ws.Range("L133").Value = 80871.375
$ws.Range("N133").Value = -90991.375
$ws.Range("H134").Value = 99995
$ws.Range("J134").Value = 99995
$ws.Range("L134").Value = 99995
$ws.Range("N134").Value = -110135
$ws.Range("H138").Value = 1449.705
$ws.Range("I138").Value = 1067.3928
$ws.Range("J138").Value = 1774.091
$ws.Range("K138").Value = 3202.1784
$ws.Range("L138").Value = 5322.272999999999
$ws.Range("M138").Value = 1937.8216
$ws.Range("N138").Value = -15602.273
$ws.Range("H139").Value = 98406.664
$ws.Range("J139").Value = 98406.664
$ws.Range("L139").Value = 98406.664
$ws.Range("N139").Value = -108686.664
$ws.Range("H140").Value = 80762
$ws.Range("J140").Value = 80762
$ws.Range("L140").Value = 80762
$ws.Range("N140").Value = -91122

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents() | Out-Null
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents() | Out-Null
$ws.Range("H104").Value = 31299
$ws.Range("J104").Value = 31299
$ws.Range("L104").Value = 31299
$ws.Range("N104").Value = -38287
$ws.Range("H121").Value = 54222.25
$ws.Range("J121").Value = 54222.25
$ws.Range("L121").Value = 54222.25
$ws.Range("N121").Value = -57716.25
$ws.Range("H131").Value = 10000
$ws.Range("J131").Value = 10000
$ws.Range("L131").Value = 10000
$ws.Range("N131").Value = -20080
$ws.Range("H135").Value = 31124.5
$ws.Range("J135").Value = 31124.5
$ws.Range("L135").Value = 31124.5
$ws.Range("N135").Value = -41264.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 43996
$ws.Range("J2").Value = 43996
$ws.Range("L2").Value = 43996
$ws.Range("N2").Value = -44222
$ws.Range("H6").Value = 6000
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents() | Out-Null
$ws.Range("H13").Value = 86158.664
$ws.Range("J13").Value = 86158.664
$ws.Range("L13").Value = 86158.664
$ws.Range("N13").Value = -86494.664
$ws.Range("H55").Value = 36994.25
$ws.Range("J55").Value = 36994.25
$ws.Range("L55").Value = 36994.25
$ws.Range("N55").Value = -37540.25
$ws.Range("H132").Value = 52304.777
$ws.Range("J132").Value = 52304.777
$ws.Range("L132").Value = 52304.777
$ws.Range("N132").Value = -62424.777
$ws.Range("H134").Value = 4622.1
$ws.Range("J134").Value = 8095
$ws.Range("L134").Value = 24285
$ws.Range("N134").Value = -29355
$ws.Range("H135").Value = 112499.664
$ws.Range("J135").Value = 112499.664
$ws.Range("L135").Value = 112499.664
$ws.Range("N135").Value = -122639.664
$ws.Range("H140").Value = 58166.484
$ws.Range("J140").Value = 43481.215
$ws.Range("L140").Value = 43481.215
$ws.Range("N140").Value = -53841.215

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 28892
$ws.Range("J18").Value = 28892
$ws.Range("L18").Value = 28892
$ws.Range("N18").Value = -29352
$ws.Range("H31").Value = 2979.0557
$ws.Range("I31").Value = 2663.4285
$ws.Range("J31").Value = 3179.9092
$ws.Range("K31").Value = 2663.4285
$ws.Range("L31").Value = 3179.9092
$ws.Range("M31").Value = -2368.4285
$ws.Range("N31").Value = -3769.9092
$ws.Range("H34").Value = 2979.0557
$ws.Range("I34").Value = 2663.4285
$ws.Range("J34").Value = 3179.9092
$ws.Range("K34").Value = 2663.4285
$ws.Range("L34").Value = 3179.9092
$ws.Range("M34").Value = -2461.4285
$ws.Range("N34").Value = -3583.9092
$ws.Range("H108").Value = 45234.855
$ws.Range("J108").Value = 45234.855
$ws.Range("L108").Value = 45234.855
$ws.Range("N108").Value = -52914.855
$ws.Range("H119").Value = 99999
$ws.Range("J119").Value = 99999
$ws.Range("L119").Value = 99999
$ws.Range("N119").Value = -109675
$ws.Range("H138").Value = 66785.5
$ws.Range("J138").Value = 66785.5
$ws.Range("L138").Value = 66785.5
$ws.Range("N138").Value = -77065.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 21.2
$ws.Range("I2").Value = 7.818182
$ws.Range("K2").Value = 46.909092
$ws.Range("M2").Value = 66.090908
$ws.Range("H13").Value = 334.33334
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 501
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1503
$ws.Range("M13").Value = 165
$ws.Range("N13").Value = -1839
$ws.Range("H134").Value = 1341.9166
$ws.Range("I134").Value = 1009.36365
$ws.Range("K134").Value = 3028.09095
$ws.Range("M134").Value = 2041.90905
$ws.Range("H139").Value = 3392.8462
$ws.Range("I139").Value = 2425.5833
$ws.Range("J139").Value = 15000
$ws.Range("K139").Value = 7276.749899999999
$ws.Range("L139").Value = 45000
$ws.Range("M139").Value = -2136.749899999999
$ws.Range("N139").Value = -55280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 10666.667
$ws.Range("J92").Value = 10666.667
$ws.Range("L92").Value = 10666.667
$ws.Range("N92").Value = -14410.667
$ws.Range("H107").Value = 839.4737
$ws.Range("I107").Value = 767.5714
$ws.Range("J107").Value = 881.4167
$ws.Range("K107").Value = 767.5714
$ws.Range("L107").Value = 881.4167
$ws.Range("M107").Value = 1152.4286
$ws.Range("N107").Value = -4721.4167
$ws.Range("H109").Value = 30847.285
$ws.Range("J109").Value = 30847.285
$ws.Range("L109").Value = 30847.285
$ws.Range("N109").Value = -32927.285
$ws.Range("H114").Value = 70557.63
$ws.Range("J114").Value = 70557.63
$ws.Range("L114").Value = 70557.63
$ws.Range("N114").Value = -79235.63
$ws.Range("H119").Value = 55305.918
$ws.Range("J119").Value = 55543.637
$ws.Range("L119").Value = 55543.637
$ws.Range("N119").Value = -65219.637

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 23166.666
$ws.Range("J10").Value = 23166.666
$ws.Range("L10").Value = 23166.666
$ws.Range("N10").Value = -23446.666
$ws.Range("H22").Value = 969.56525
$ws.Range("I22").Value = 931
$ws.Range("K22").Value = 931
$ws.Range("M22").Value = -636
$ws.Range("H27").Value = 969.56525
$ws.Range("I27").Value = 931
$ws.Range("K27").Value = 931
$ws.Range("M27").Value = -824
$ws.Range("H40").Value = 7411065.5
$ws.Range("I40").Value = 4113.3125
$ws.Range("K40").Value = 4113.3125
$ws.Range("M40").Value = -3977.3125
$ws.Range("H55").Value = 1018.93335
$ws.Range("I55").Value = 640.9091
$ws.Range("J55").Value = 2058.5
$ws.Range("K55").Value = 640.9091
$ws.Range("L55").Value = 2058.5
$ws.Range("M55").Value = -467.9091
$ws.Range("N55").Value = -2404.5
$ws.Range("H129").Value = 93376.55499999999
$ws.Range("J129").Value = 101999.86
$ws.Range("L129").Value = 101999.86
$ws.Range("N129").Value = -111999.86

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 12666
$ws.Range("J20").Value = 12666
$ws.Range("L20").Value = 12666
$ws.Range("N20").Value = -13146
$ws.Range("H22").Value = 14749.5
$ws.Range("J22").Value = 14749.5
$ws.Range("L22").Value = 14749.5
$ws.Range("N22").Value = -15335.5
$ws.Range("H122").Value = 1729.5416
$ws.Range("I122").Value = 1358.3684
$ws.Range("J122").Value = 3140
$ws.Range("K122").Value = 4075.1052
$ws.Range("L122").Value = 9420
$ws.Range("M122").Value = -1625.1052
$ws.Range("N122").Value = -14320
$ws.Range("H136").Value = 1759.7333
$ws.Range("I136").Value = 1603.2727
$ws.Range("J136").Value = 2190
$ws.Range("K136").Value = 4809.8181
$ws.Range("L136").Value = 6570
$ws.Range("M136").Value = -2259.8181
$ws.Range("N136").Value = -11670
